$wb = $excel.ActiveWorkbook

# --- Select the full data range on Italy before we branch off Spain, so
#     that Italy ends up with a non-active, full-range selection (as in
#     the target) once Spain becomes the new active sheet. ---
$italy = $wb.Worksheets.Item("Italy")
$italy.Activate()
$italy.Range("A1:D10").Select()

# --- Create the new "Spain" sheet by copying "Italy" (keeps styles,
#     merged cells, column layout, etc.) and placing it right after it. ---
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# --- Market-specific content ---
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2064/T2063"

# --- Column widths (best-fit-like values) ---
$spain.Columns.Item(1).ColumnWidth = 25.109375
$spain.Columns.Item(2).ColumnWidth = 36.44140625
$spain.Columns.Item(4).ColumnWidth = 19.33203125

# --- Row heights for the wrapped rows ---
$spain.Rows.Item(3).RowHeight = 28.8
$spain.Rows.Item(4).RowHeight = 28.8
$spain.Rows.Item(5).RowHeight = 28.8

# --- Make Spain the active sheet/tab with its own lingering selection ---
$spain.Activate()
$spain.Range("C14").Select()
